$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.965.39"
$ws.Range("E2").Value = "  +4.40%  "
# Row 3
$ws.Range("D3").Value = "2.285.37"
$ws.Range("E3").Value = "  +5.20%  "
# Row 4
$ws.Range("E4").Value = "  +0.10%  "
# Row 5
$ws.Range("E5").Value = "  +0.40%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +3.27%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.49"
$ws.Range("E7").Value = "  +9.98%  "
# Row 8
$ws.Range("E8").Value = "  -0.07%  "
# Row 9
$ws.Range("E9").Value = "  +15.20%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("E10").Value = "  +7.92%  "
# Row 11
$ws.Range("E11").Value = "  +5.33%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.85"
$ws.Range("E12").Value = "  +1.53%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.40"
$ws.Range("E13").Value = "  +8.60%  "
# Row 14
$ws.Range("E14").Value = "  +3.02%  "
# Row 15
$ws.Range("D15").Value = "2.625.59"
$ws.Range("E15").Value = "  +5.14%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.11"
$ws.Range("E16").Value = "  +6.22%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("E17").Value = "  +5.72%  "
# Row 18
$ws.Range("D18").Value = "2.285.82"
$ws.Range("E18").Value = "  +5.53%  "
# Row 19
$ws.Range("D19").Value = "42.908.84"
$ws.Range("E19").Value = "  +4.56%  "
# Row 20
$ws.Range("E20").Value = "  +7.58%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  +6.04%  "
# Row 22
$ws.Range("E22").Value = "  +2.90%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.86"
$ws.Range("E23").Value = "  +2.75%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +6.92%  "
# Row 25
$ws.Range("E25").Value = "  +1.89%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.69"
$ws.Range("E26").Value = "  +3.15%  "
# Row 27
$ws.Range("E27").Value = "  -0.01%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.46"
$ws.Range("E28").Value = "  +2.03%  "
# Row 29
$ws.Range("E29").Value = "  -0.90%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  +1.15%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.59"
$ws.Range("E31").Value = "  -0.08%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.15"
$ws.Range("E32").Value = "  +5.13%  "
# Row 33
$ws.Range("E33").Value = "  +11.36%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("E34").Value = "  +7.06%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0814"
$ws.Range("E35").Value = "  +9.32%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.27"
$ws.Range("E36").Value = "  +28.27%  "
# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.80"
$ws.Range("E37").Value = "  +22.89%  "
# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.126"
$ws.Range("E38").Value = "  +4.33%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.81"
$ws.Range("E39").Value = "  +6.54%  "
# Row 40
$ws.Range("E40").Value = "  +2.17%  "
# Row 41
$ws.Range("E41").Value = "  +6.06%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.24"
$ws.Range("E42").Value = "  +16.88%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.05"
$ws.Range("E43").Value = "  +10.57%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.215"
$ws.Range("E44").Value = "  +15.47%  "
# Row 45
$ws.Range("E45").Value = "  +9.27%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.00"
$ws.Range("E46").Value = "  -9.26%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.82"
$ws.Range("E47").Value = "  +1.70%  "
# Row 48
$ws.Range("E48").Value = "  +5.63%  "
# Row 49
$ws.Range("E49").Value = "  +5.88%  "
# Row 50
$ws.Range("E50").Value = "  -0.04%  "
# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("E51").Value = "  +5.44%  "
